$wb = $excel.ActiveWorkbook

# "Forecast Comparison" sheet: update MyForecast (column D) values for weeks 12-16 (rows 13-17)
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsForecast.Range("D13").Value = 68
$wsForecast.Range("D14").Value = 84
$wsForecast.Range("D15").Value = 91
$wsForecast.Range("D16").Value = 82
$wsForecast.Range("D17").Value = 61

# "Summary" sheet: update Total Forecast (16 Weeks) and Min Forecast to reflect new values.
# These cells store numbers-as-text, so a leading apostrophe keeps them as text
# (matching the workbook's existing inline-string/text cell type) instead of
# letting Excel auto-convert them to numeric values.
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B9").Value = "'1950"
$wsSummary.Range("B14").Value = "'61"
